$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New arrival rows (Sunday, Jan 15) appended after the existing data, rows 76-81.
# Column layout matches existing sheet: A=NUMBER, B=DATE, C=TIME, D=FLIGHT, E=FROM,
# F=SHORT, G=AIRLINE, H=MODEL, I=AIRCFAT ID, J=STATUS, K=(blank), L=DIFFERENCE, M=(blank)

$ws.Cells.Item(76, 1).Value = 75
$ws.Cells.Item(76, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(76, 3).Value = "6:26 AM"
$ws.Cells.Item(76, 4).Value = "UNKNOWN"
$ws.Cells.Item(76, 5).Value = "Skopje"
$ws.Cells.Item(76, 6).Value = "(SKP)"
$ws.Cells.Item(76, 7).Value = "Wizz Air "
$ws.Cells.Item(76, 8).Value = "A320"
$ws.Cells.Item(76, 9).Value = "(HA-LYG)"
$ws.Cells.Item(76, 10).Value = "6:48 AM"
$ws.Cells.Item(76, 12).Value = "0 hours, 22 minutes"

$ws.Cells.Item(77, 1).Value = 76
$ws.Cells.Item(77, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(77, 3).Value = "9:47 AM"
$ws.Cells.Item(77, 4).Value = "5Y317"
$ws.Cells.Item(77, 5).Value = "Fairfield"
$ws.Cells.Item(77, 6).Value = "(SUU)"
$ws.Cells.Item(77, 7).Value = "Atlas Air "
$ws.Cells.Item(77, 8).Value = "B744"
$ws.Cells.Item(77, 9).Value = "(N485MC)"
$ws.Cells.Item(77, 10).Value = "11:02 AM"
$ws.Cells.Item(77, 12).Value = "1 hours, 15 minutes"

$ws.Cells.Item(78, 1).Value = 77
$ws.Cells.Item(78, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(78, 3).Value = "11:20 AM"
$ws.Cells.Item(78, 4).Value = "FR2136"
$ws.Cells.Item(78, 5).Value = "London"
$ws.Cells.Item(78, 6).Value = "(STN)"
$ws.Cells.Item(78, 7).Value = "Lauda Europe "
$ws.Cells.Item(78, 8).Value = "A320"
$ws.Cells.Item(78, 9).Value = "(9H-LOT)"
$ws.Cells.Item(78, 10).Value = "11:00 AM"
$ws.Cells.Item(78, 12).Value = "0 hours, -20 minutes"

$ws.Cells.Item(79, 1).Value = 78
$ws.Cells.Item(79, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(79, 3).Value = "2:15 PM"
$ws.Cells.Item(79, 4).Value = "LO3801"
$ws.Cells.Item(79, 5).Value = "Warsaw"
$ws.Cells.Item(79, 6).Value = "(WAW)"
$ws.Cells.Item(79, 7).Value = "LOT "
$ws.Cells.Item(79, 8).Value = "E190"
$ws.Cells.Item(79, 9).Value = "(SP-LMA)"
$ws.Cells.Item(79, 10).Value = "2:10 PM"
$ws.Cells.Item(79, 12).Value = "0 hours, -5 minutes"

$ws.Cells.Item(80, 1).Value = 79
$ws.Cells.Item(80, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(80, 3).Value = "3:01 PM"
$ws.Cells.Item(80, 4).Value = "VMP944"
$ws.Cells.Item(80, 5).Value = "Paris"
$ws.Cells.Item(80, 6).Value = "(LBG)"
$ws.Cells.Item(80, 7).Value = "Execujet Scandinavia "
$ws.Cells.Item(80, 8).Value = "GL5T"
$ws.Cells.Item(80, 9).Value = "(OY-VIZ)"
$ws.Cells.Item(80, 10).Value = "3:04 PM"
$ws.Cells.Item(80, 12).Value = "0 hours, 3 minutes"

$ws.Cells.Item(81, 1).Value = 80
$ws.Cells.Item(81, 2).Value = "Sunday, Jan 15"
$ws.Cells.Item(81, 3).Value = "4:05 PM"
$ws.Cells.Item(81, 4).Value = "LO3807"
$ws.Cells.Item(81, 5).Value = "Warsaw"
$ws.Cells.Item(81, 6).Value = "(WAW)"
$ws.Cells.Item(81, 7).Value = "LOT (Star Alliance Livery) "
$ws.Cells.Item(81, 8).Value = "E170"
$ws.Cells.Item(81, 9).Value = "(SP-LDK)"
$ws.Cells.Item(81, 10).Value = "4:06 PM"
$ws.Cells.Item(81, 12).Value = "0 hours, 1 minutes"
